# updated sequence diagram for choose command
#
# Applies the diff to LogicComponentSequenceDiagram.pptx (3rd slide):
#  1. Move "Rectangle 72" (lifeline activation box).
#  2. Move + recolor + re-dash "Straight Arrow Connector 75".
#  3. Resize/reposition "TextBox 77" and split its text into 3 runs
#     ("post(" + "JumpToBrowserListEvent" + "))").
#  4. Reposition/resize "Rectangle 62" (result:Command box, id 84).
#  5. Move "Straight Arrow Connector 12".
#  6. Move "Rectangle 2" (activation box, id 3).
#  7. Add two new dashed/dash-dot arrow connectors.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------
# 1) Rectangle 72 (shape #23, id 73) - small offset nudge
# ---------------------------------------------------------------------
$rect72 = $s.Shapes.Item(23)
$rect72.Left = 811.0435
$rect72.Top  = 356.0067

# ---------------------------------------------------------------------
# 2) Straight Arrow Connector 75 (shape #25, id 76)
#    move up, recolor accent4(60/40) -> accent3(75), dash sysDash -> solid
# ---------------------------------------------------------------------
$conn75 = $s.Shapes.Item(25)
$conn75.Left = 702
$conn75.Top  = 356.86962890625
$conn75.Line.ForeColor.ObjectThemeColor = 7   # accent3
$conn75.Line.DashStyle = 1                    # solid

# ---------------------------------------------------------------------
# 3) TextBox 77 (shape #27, id 78) - reposition/resize + new text runs
# ---------------------------------------------------------------------
$tb77 = $s.Shapes.Item(27)
$tb77.Left   = 522
$tb77.Top    = 336
$tb77.Width  = 176.70591
$tb77.Height = 33.9281102

$tr = $tb77.TextFrame.TextRange
$tr.Text = "post("
$tr.InsertAfter("JumpToBrowserListEvent") | Out-Null
$tr.InsertAfter("))") | Out-Null

# ---------------------------------------------------------------------
# 4) Rectangle 62 / "result:Command" box (shape #32, id 84)
# ---------------------------------------------------------------------
$rect84 = $s.Shapes.Item(32)
$rect84.Left   = 528.7752
$rect84.Top    = 397.0435
$rect84.Width  = 168.2152099609375
$rect84.Height = 16.6417323

# ---------------------------------------------------------------------
# 5) Straight Arrow Connector 12 (shape #42, id 13)
# ---------------------------------------------------------------------
$conn12 = $s.Shapes.Item(42)
$conn12.Left = 516
$conn12.Top  = 402

# ---------------------------------------------------------------------
# 6) Rectangle 2 (shape #47, id 3) - activation bar
# ---------------------------------------------------------------------
$rect2 = $s.Shapes.Item(47)
$rect2.Left = 696
$rect2.Top  = 354.9786

# ---------------------------------------------------------------------
# 7) New connector: "Straight Arrow Connector 53"
#    accent4 (lumMod60/lumOff40), dashDot
# ---------------------------------------------------------------------
$newConn1 = $s.Shapes.AddConnector(1, 0, 0, 100, 0)
$newConn1.Name = "Straight Arrow Connector 53"
$newConn1.Left   = 710.8696
$newConn1.Top    = 371.21734619140625
$newConn1.Width  = 111.1304
$newConn1.Height = 0

$ln1 = $newConn1.Line
$ln1.Weight = 1.5
$ln1.ForeColor.ObjectThemeColor = 8   # accent4
$ln1.DashStyle = 5                    # dashDot
$ln1.BeginArrowheadStyle = 3          # arrow
$ln1.BeginArrowheadLength = 2         # med
$ln1.BeginArrowheadWidth  = 2         # med
$ln1.EndArrowheadStyle = 1            # none
$ln1.EndArrowheadLength = 2
$ln1.EndArrowheadWidth  = 2

# ---------------------------------------------------------------------
# 8) New connector: "Straight Arrow Connector 55"
#    accent3 (lumMod75), dash
# ---------------------------------------------------------------------
$newConn2 = $s.Shapes.AddConnector(1, 0, 0, 100, 0)
$newConn2.Name = "Straight Arrow Connector 55"
$newConn2.Left   = 515.27354
$newConn2.Top    = 376.95654
$newConn2.Width  = 192.7265
$newConn2.Height = 0

$ln2 = $newConn2.Line
$ln2.Weight = 1.5
$ln2.ForeColor.ObjectThemeColor = 7   # accent3
$ln2.DashStyle = 4                    # dash
$ln2.BeginArrowheadStyle = 3          # arrow
$ln2.BeginArrowheadLength = 2         # med
$ln2.BeginArrowheadWidth  = 2         # med
$ln2.EndArrowheadStyle = 1            # none
$ln2.EndArrowheadLength = 2
$ln2.EndArrowheadWidth  = 2
